$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Gracz 3"
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Gracz 5"
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Gracz test"
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Gracz 1"
$ws.Range("B6").Value = "pkstz"
$ws.Range("D6").Value = "Zachodniopomorskie"
$ws.Range("E6").Value = 95
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Gracz 11"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Gracz a"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Gracz 8"
$ws.Range("D9").Value = "Opolskie"
$ws.Range("E9").Value = 92
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Gracz 26"
$ws.Range("D10").Value = "Podlaskie"
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Gracz test 2"
$ws.Range("D11").Value = "Lubelskie"
$ws.Range("E11").Value = 88
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Maks"
$ws.Range("D12").Value = "Dolnośląskie"
$ws.Range("E12").Value = 87
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "maks"
$ws.Range("D13").Value = "Wszystkie"
$ws.Range("E13").Value = 85
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "k"
$ws.Range("D14").Value = "Opolskie"
$ws.Range("E14").Value = 83
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Maks"
$ws.Range("C15").Value = "Extreme"
$ws.Range("D15").Value = "Lubelskie"
$ws.Range("E15").Value = 72
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "mak"
$ws.Range("C16").Value = "Easy"
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "l"
$ws.Range("D17").Value = "Opolskie"
$ws.Range("E17").Value = 67
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "pkstz"
$ws.Range("D18").Value = "Zachodniopomorskie"
$ws.Range("E18").Value = 66
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "wiki <3"
$ws.Range("D19").Value = "Lubelskie"
$ws.Range("E19").Value = 60
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "testyyyyyyy"
$ws.Range("C20").Value = "Extreme"
$ws.Range("D20").Value = "Kujawsko-Pomorskie"
$ws.Range("E20").Value = 57
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "k"
$ws.Range("C21").Value = "Hard"
$ws.Range("D21").Value = "Lubuskie"
$ws.Range("E21").Value = 50
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Maks"
$ws.Range("D22").Value = "Podlaskie"
$ws.Range("E22").Value = 48
$ws.Range("A23").Value = 24
$ws.Range("B25").Value = "a"
$ws.Range("D25").Value = "Lubelskie"
$ws.Range("E25").Value = 42
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "maks"
$ws.Range("D26").Value = "Wszystkie"
$ws.Range("E26").Value = 40
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = "Extreme"
$ws.Range("D28").Value = "Dolnośląskie"
$ws.Range("B29").Value = "a"
$ws.Range("C29").Value = "Hard"
$ws.Range("D29").Value = "Lubuskie"
$ws.Range("E29").Value = 36
$ws.Range("B30").Value = "kkk"
$ws.Range("C30").Value = "Easy"
$ws.Range("D30").Value = "Podlaskie"
$ws.Range("B31").Value = "Maks"
$ws.Range("C31").Value = "Medium"
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = "Extreme"
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "szymek"
$ws.Range("D33").Value = "Kujawsko-Pomorskie"
$ws.Range("E33").Value = 35
$ws.Range("B34").Value = "pkstz"
$ws.Range("C34").Value = "Extreme"
$ws.Range("D34").Value = "Zachodniopomorskie"
$ws.Range("A35").Value = 41
$ws.Range("B35").Value = "h"
$ws.Range("D35").Value = "Opolskie"
$ws.Range("A36").Value = 40
$ws.Range("B36").Value = "maks"
$ws.Range("D36").Value = "Śląskie"
$ws.Range("A37").Value = 39
$ws.Range("B37").Value = "hjk"
$ws.Range("D37").Value = "Lubelskie"
$ws.Range("A38").Value = 38
$ws.Range("B38").Value = ""
$ws.Range("D38").Value = "Opolskie"
$ws.Range("A39").Value = 0
$ws.Range("B39").Value = "k"
$ws.Range("C39").Value = "Extreme"
$ws.Range("D39").Value = "Małopolskie"
$ws.Range("A40").Value = 36
$ws.Range("B40").Value = "uj"
$ws.Range("C40").Value = "Medium"
$ws.Range("D40").Value = "Opolskie"
$ws.Range("A41").Value = 35
$ws.Range("B41").Value = "d"
$ws.Range("D41").Value = "Lubelskie"
$ws.Range("E41").Value = 33
$ws.Range("A42").Value = 34
$ws.Range("B42").Value = "maks"
$ws.Range("D42").Value = "Śląskie"
$ws.Range("E42").Value = 33
$ws.Range("A43").Value = 33
$ws.Range("B43").Value = "h"
$ws.Range("C43").Value = "Hard"
$ws.Range("D43").Value = "Małopolskie"
$ws.Range("E43").Value = 33
$ws.Range("B44").Value = "aa"
$ws.Range("D44").Value = "Dolnośląskie"
$ws.Range("E44").Value = 32
$ws.Range("B45").Value = "''"
$ws.Range("D45").Value = "Łódzkie"
$ws.Range("E45").Value = 32
$ws.Range("C46").Value = "Extreme"
$ws.Range("D46").Value = "Dolnośląskie"
$ws.Range("E46").Value = 30
$ws.Range("B47").Value = "k"
$ws.Range("D47").Value = "Lubuskie"
$ws.Range("E47").Value = 29
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = "Medium"
$ws.Range("D48").Value = "Lubuskie"
$ws.Range("E48").Value = 29
$ws.Range("B49").Value = "hjk"
$ws.Range("D49").Value = "Lubelskie"
$ws.Range("E49").Value = 29
$ws.Range("B50").Value = "Maks"
$ws.Range("D50").Value = "Wszystkie"
$ws.Range("E50").Value = 28
$ws.Range("B51").Value = "iu"
$ws.Range("D51").Value = "Łódzkie"
$ws.Range("E51").Value = 28
